$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the (now orphaned) "_GoBack" bookmark that currently sits
#    in its own empty paragraph right after the "9. Sys" paragraph.
#    We locate that paragraph by text rather than by a hard-coded
#    index so the script is resilient to other structural details.
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$sysPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "9. Sys*") {
        $sysPara = $p
        break
    }
}

if ($sysPara -ne $null) {
    $bookmarkPara = $sysPara.Next()
    # Deleting the range clears the paragraph's contents (the
    # bookmarkStart/bookmarkEnd pair) while leaving the empty
    # paragraph itself in place.
    $bookmarkPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2. In the Alternative-Flow / A1 paragraph, insert " form" right
#    after " a valid email address" (before the closing curly quote),
#    then re-insert the "_GoBack" bookmark right after the newly
#    added text.
# ------------------------------------------------------------------
$findRng = $d.Content.Duplicate
$found = $findRng.Find.Execute("a valid email address", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPos = $findRng.End
    $insertRng = $d.Range($insertPos, $insertPos)
    $insertRng.InsertBefore(" form")

    # Force the newly inserted text into its own run (matching the
    # plain, non-bold formatting of its neighbours) instead of being
    # silently merged into the preceding run.
    $newTextRng = $d.Range($insertPos, $insertPos + 5)
    $newTextRng.Bold = 1
    $newTextRng.Bold = 0

    $bmPos = $insertPos + 5
    $bmRng = $d.Range($bmPos, $bmPos)
    $bmRng.Bookmarks.Add("_GoBack")
}
